# NYPD CompStat weekly report refresh — new crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/number and reporting week dates ---
$ws.Range("A8").Value = "Volume 31   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/1/2024  Through  7/7/2024"

# --- Row 14: Murder ---
$ws.Range("N14").Value = -93.333333333333

# --- Row 15: Rape ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0

# --- Row 16: Robbery ---
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -18.75
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 101
$ws.Range("K16").Value = -4.950495049504
$ws.Range("L16").Value = 3.225806451612
$ws.Range("M16").Value = -53.170731707317
$ws.Range("N16").Value = -90.742526518804

# --- Row 17: Fel. Assault ---
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -61.538461538461
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = -26.190476190476
$ws.Range("I17").Value = 219
$ws.Range("J17").Value = 223
$ws.Range("K17").Value = -1.793721973094
$ws.Range("L17").Value = -7.594936708860
$ws.Range("M17").Value = -6.008583690987
$ws.Range("N17").Value = -62.371134020618

# --- Row 18: Burglary ---
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 95
$ws.Range("J18").Value = 115
$ws.Range("K18").Value = -17.391304347826
$ws.Range("L18").Value = -19.491525423728
$ws.Range("M18").Value = -57.399103139013
$ws.Range("N18").Value = -85.271317829457

# --- Row 19: Gr. Larceny ---
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 30
$ws.Range("H19").Value = -11.764705882352
$ws.Range("I19").Value = 189
$ws.Range("J19").Value = 218
$ws.Range("K19").Value = -13.302752293578
$ws.Range("L19").Value = -2.577319587628
$ws.Range("M19").Value = -12.5
$ws.Range("N19").Value = -52.512562814070

# --- Row 20: G.L.A. ---
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 62.5
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 62
$ws.Range("K20").Value = -16.129032258064
$ws.Range("L20").Value = 13.043478260869
$ws.Range("M20").Value = -24.637681159420
$ws.Range("N20").Value = -81.944444444444

# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -44.736842105263
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 124
$ws.Range("H21").Value = -19.354838709677
$ws.Range("I21").Value = 665
$ws.Range("J21").Value = 737
$ws.Range("K21").Value = -9.769335142469
$ws.Range("L21").Value = -6.338028169014
$ws.Range("M21").Value = -31.088082901554
$ws.Range("N21").Value = -77.958236658932

# --- Row 22: Transit ---
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = -50
$ws.Range("L22").Value = -64.285714285714
$ws.Range("M22").Value = -79.166666666666

# --- Row 23: Housing ---
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = -38.888888888888
$ws.Range("I23").Value = 123
$ws.Range("J23").Value = 127
$ws.Range("K23").Value = -3.149606299212
$ws.Range("L23").Value = 1.652892561983
$ws.Range("M23").Value = 11.818181818181

# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 7.692307692307
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = -2.564102564102
$ws.Range("I24").Value = 688
$ws.Range("J24").Value = 773
$ws.Range("K24").Value = -10.996119016817
$ws.Range("L24").Value = -12.244897959183
$ws.Range("M24").Value = 43.035343035343

# --- Row 25: Retail Theft ---
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = 26.086956521739
$ws.Range("I25").Value = 281
$ws.Range("J25").Value = 302
$ws.Range("K25").Value = -6.953642384105
$ws.Range("L25").Value = -18.786127167630

# --- Row 26: Misd. Assault ---
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 128.571428571429
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = 39.024390243902
$ws.Range("I26").Value = 342
$ws.Range("J26").Value = 328
$ws.Range("K26").Value = 4.268292682926
$ws.Range("L26").Value = -6.811989100817
$ws.Range("M26").Value = -31.048387096774

# --- Row 27: UCR Rape* ---
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 17
$ws.Range("K27").Value = -19.047619047619
$ws.Range("L27").Value = -34.615384615384

# --- Row 28: Other Sex Crimes ---
# C28 goes from numeric 2 to the text placeholder "0" (shared with C14/C22/etc.);
# copy it in from an existing placeholder cell so style + shared string are reused exactly.
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = -20
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = 33.333333333333

# --- Row 29: Shooting Vic. ---
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("N29").Value = -91.156462585034

# --- Row 30: Shooting Inc. ---
$ws.Range("C14").Copy($ws.Range("C30"))
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -75
$ws.Range("N30").Value = -91.338582677165
